# Refresh the Price (D) / Volume(1h) (E) columns of the cryptos sheet with
# the latest scrape, per commit "Updated cryptos list on Fri Jul  5 20:49:09
# UTC 2024 with GitHub Actions". Only the cell text changes; everything else
# (rank, coin name, link, layout, styling) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.344.05"
$ws.Range("E2").Value = "  -3.17%  "

$ws.Range("D3").Value = "2.967.79"
$ws.Range("E3").Value = "  -5.34%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "495.08"
$ws.Range("E5").Value = "  -5.61%  "

$ws.Range("D6").Value = "134.42"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "2.968.03"
$ws.Range("E8").Value = "  -5.29%  "

$ws.Range("E9").Value = "  -4.08%  "

$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("E11").Value = "  -3.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.350"
$ws.Range("E12").Value = "  -7.24%  "

$ws.Range("E13").Value = "  -0.74%  "

$ws.Range("D14").Value = "3.475.47"
$ws.Range("E14").Value = "  -5.33%  "

$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("D16").Value = "56.326.29"
$ws.Range("E16").Value = "  -3.19%  "

$ws.Range("D17").Value = "2.966.98"
$ws.Range("E17").Value = "  -5.32%  "

$ws.Range("D18").Value = "0.0000146"
$ws.Range("E18").Value = "  -4.34%  "

$ws.Range("D19").Value = "5.77"
$ws.Range("E19").Value = "  +0.67%  "

$ws.Range("D20").Value = "12.33"
$ws.Range("E20").Value = "  -5.16%  "

$ws.Range("E21").Value = "  -1.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "325.50"
$ws.Range("E22").Value = "  -5.24%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").Value = "0.467"
$ws.Range("E24").Value = "  -7.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "61.40"
$ws.Range("E25").Value = "  -9.47%  "

$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.45%  "

$ws.Range("D27").Value = "0.161"
$ws.Range("E27").Value = "  -5.86%  "

$ws.Range("D28").Value = "0.0₃0895"
$ws.Range("E28").Value = "  -5.85%  "

$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").Value = "6.47"
$ws.Range("E30").Value = "  -5.21%  "

$ws.Range("E31").Value = "  -2.78%  "

$ws.Range("E32").Value = "  -5.81%  "

$ws.Range("E33").Value = "  -6.85%  "

$ws.Range("E34").Value = "  -5.66%  "

$ws.Range("D35").Value = "152.56"
$ws.Range("E35").Value = "  -3.04%  "

$ws.Range("E36").Value = "  -8.32%  "

$ws.Range("E37").Value = "  -7.07%  "

$ws.Range("E38").Value = "  -10.39%  "

$ws.Range("E39").Value = "  -2.42%  "

$ws.Range("E40").Value = "  -2.72%  "

$ws.Range("D41").Value = "3.001.02"
$ws.Range("E41").Value = "  -5.19%  "

$ws.Range("D42").Value = "36.58"
$ws.Range("E42").Value = "  -9.55%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("E44").Value = "  -7.73%  "

$ws.Range("D45").Value = "0.991"
$ws.Range("E45").Value = "  -9.09%  "

$ws.Range("D46").Value = "2.201.37"
$ws.Range("E46").Value = "  -3.54%  "

$ws.Range("E47").Value = "  -3.58%  "

$ws.Range("E48").Value = "  -9.08%  "

$ws.Range("D49").Value = "1.93"
$ws.Range("E49").Value = "  +3.80%  "

$ws.Range("E50").Value = "  +1.03%  "

$ws.Range("E51").Value = "  -6.89%  "
